$wb = $excel.ActiveWorkbook

# Sheet references (by name, matches xl/workbook.xml sheet order):
#   1 = lower, 2 = external, 3 = internal, 4 = upper, 5 = final ini
$wsLower    = $wb.Worksheets.Item("lower")
$wsExternal = $wb.Worksheets.Item("external")
$wsInternal = $wb.Worksheets.Item("internal")
$wsUpper    = $wb.Worksheets.Item("upper")
$wsFinal    = $wb.Worksheets.Item("final ini")

# --- Data / input-cell edits -------------------------------------------------
# These drive the downstream formulas (C/D/E/G columns on lower/external/
# internal/upper, and the fully-derived "final ini" sheet) via automatic
# recalculation, so only the raw inputs need to be written.

# lower: flip the R3 sign convention and move the N7 offset
$wsLower.Range("R3").Value = -1
$wsLower.Range("N7").Value = -88

# external: flip R3, move N7 and N9
$wsExternal.Range("R3").Value = -1
$wsExternal.Range("N7").Value = 53
$wsExternal.Range("N9").Value = 300

# internal: flip R3, move N7 and N9
$wsInternal.Range("R3").Value = -1
$wsInternal.Range("N7").Value = 90
$wsInternal.Range("N9").Value = -52

# upper: move N7 and N9 (R3 stays at 1)
$wsUpper.Range("N7").Value = -186
$wsUpper.Range("N9").Value = -110

# --- View / selection edits --------------------------------------------------
# Apply per-sheet zoom + selection, then activate sheets in order so the
# last Activate() call leaves "lower" as the selected tab (matches the
# workbook bookView's activeTab becoming sheet index 0).

$wsExternal.Activate()
$excel.ActiveWindow.Zoom = 55
$wsExternal.Range("N10").Select()

$wsInternal.Activate()
$excel.ActiveWindow.Zoom = 55
$wsInternal.Range("N10").Select()

$wsUpper.Activate()
$excel.ActiveWindow.Zoom = 85
$wsUpper.Range("N9").Select()

$wsFinal.Activate()
$excel.ActiveWindow.Zoom = 85
$wsFinal.Range("A3:G94").Select()

$wsLower.Activate()
$excel.ActiveWindow.Zoom = 55
$wsLower.Range("N7").Select()
